$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to Text format so numeric-looking strings are preserved as text
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "23.076.00"
$ws.Range("E2").Value = "  -3.64%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.598.78"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "300.78"
$ws.Range("E6").Value = "  -2.91%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.3773"
$ws.Range("E7").Value = "  -2.76%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3629"
$ws.Range("E8").Value = "  -4.48%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "49.91"
$ws.Range("E9").Value = "  -4.36%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.249"
$ws.Range("E10").Value = "  -6.79%  "

$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.08111"
$ws.Range("E12").Value = "  -3.86%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "22.91"
$ws.Range("E13").Value = "  -3.92%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.550"
$ws.Range("E14").Value = "  -6.99%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.00001253"
$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.337"
$ws.Range("E16").Value = "  -8.43%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.600.94"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "91.55"
$ws.Range("E18").Value = "  -2.71%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.06874"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "18.22"
$ws.Range("E20").Value = "  -6.98%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.534"
$ws.Range("E21").Value = "  -5.77%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "12.93"
$ws.Range("E23").Value = "  -5.82%  "

$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "23.090.12"
$ws.Range("E24").Value = "  -3.60%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.345"
$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.707"
$ws.Range("E26").Value = "  -7.99%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.00"
$ws.Range("E27").Value = "  -4.62%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "149.80"
$ws.Range("E28").Value = "  -2.34%  "

$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").Value = "5.245"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "131.64"
$ws.Range("E30").Value = "  -4.55%  "

$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "2.423"
$ws.Range("E31").Value = "  -3.47%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "6.791"
$ws.Range("E32").Value = "  -13.11%  "

$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.777.64"
$ws.Range("E33").Value = "  -2.76%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.9501"
$ws.Range("E34").Value = "  -6.53%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.07632"
$ws.Range("E35").Value = "  -4.93%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02716"
$ws.Range("E36").Value = "  -7.19%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2536"
$ws.Range("E37").Value = "  -5.06%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "6.191"
$ws.Range("E38").Value = "  -8.17%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.08875"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "9.996"
$ws.Range("E40").Value = "  -6.43%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.367"
$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.7047"
$ws.Range("E42").Value = "  -6.85%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "12.58"
$ws.Range("E43").Value = "  -5.83%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "15.41"
$ws.Range("E44").Value = "  -5.23%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6570"
$ws.Range("E45").Value = "  -5.35%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.285"
$ws.Range("E47").Value = "  -6.33%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.977"
$ws.Range("E48").Value = "  -2.70%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.07940"
$ws.Range("E49").Value = "  -4.38%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "131.10"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.201"
$ws.Range("E51").Value = "  -2.20%  "
